# Commit: feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after "总计" (before the current
#    "2022-Q2" sheet), populated with the fund-holdings table for that quarter.
# 2. Insert a corresponding new row at the top of the "总计" (summary) sheet's
#    data table for "2022-Q3" (12 holdings, 4.36 亿元), pushing the existing
#    quarters down by one row.

$wb = $excel.ActiveWorkbook
$ws_total = $wb.Worksheets.Item(1)
$ws_q2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Step 1: shift the existing data rows (2..8) of "总计" down to (3..9)
# working from the bottom up so we don't overwrite data before reading it.
# ---------------------------------------------------------------------
for ($r = 8; $r -ge 2; $r--) {
    $rNew = $r + 1
    $ws_total.Cells.Item($rNew, 2).Value = $ws_total.Cells.Item($r, 2).Value2
    $ws_total.Cells.Item($rNew, 3).Value = $ws_total.Cells.Item($r, 3).Value2
    $ws_total.Cells.Item($rNew, 4).Value = $ws_total.Cells.Item($r, 4).Value2
    $ws_total.Cells.Item($rNew, 1).Value = $r - 1
}

# New first data row (row 2) describing the 2022-Q3 quarter.
$ws_total.Cells.Item(2, 1).Value = 0
$ws_total.Cells.Item(2, 2).Value = "2022-Q3"
$ws_total.Cells.Item(2, 3).Value = 12
$ws_total.Cells.Item(2, 4).Value = 4.36

# Row 9 is brand new (the table used to stop at row 8), so its "A" cell
# needs the bold/bordered index-column style copied across explicitly.
$ws_total.Range("A8").Copy()
$ws_total.Range("A9").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 2: create the new "2022-Q3" worksheet right before "2022-Q2" and
# fill it in with the per-fund holdings table.
# ---------------------------------------------------------------------
$ws_new = $wb.Worksheets.Add($ws_q2)
$ws_new.Name = "2022-Q3"
$ws_new.Outline.SummaryRow = 1
$ws_new.Outline.SummaryColumn = 1

# Header row
$ws_new.Cells.Item(1, 2).Value = "基金代码"
$ws_new.Cells.Item(1, 3).Value = "基金名称"
$ws_new.Cells.Item(1, 4).Value = "基金规模"
$ws_new.Cells.Item(1, 5).Value = "股票总仓位"
$ws_new.Cells.Item(1, 6).Value = "仓位占比"
$ws_new.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws_new.Cells.Item(1, 8).Value = "仓位排名"

# Apply the bold/bordered/centered header style (same style used for the
# index column, copied from the "总计" sheet's A2 cell which already has it).
$styleSrc = $ws_total.Range("A2")
$styleSrc.Copy()
$ws_new.Range("B1:H1").PasteSpecial(-4122)
$ws_new.Range("A2:A13").PasteSpecial(-4122)

# Columns B,D,E,F,G hold numeric-looking values (fund codes, sizes, ratios)
# that must nonetheless be stored as *text*, exactly like all the other
# per-quarter sheets in this workbook. Force a text number format on those
# ranges before writing so Excel doesn't silently coerce them to numbers
# (and so leading zeros in fund codes like "014509" survive).
$ws_new.Range("B2:B13").NumberFormat = "@"
$ws_new.Range("D2:G12").NumberFormat = "@"
$ws_new.Range("D13:F13").NumberFormat = "@"

$rows = @(
    @(0,  "570001", "诺德价值优势混合",             "42.75", "92.27", "5.47", "2.3384", 10),
    @(1,  "015182", "汇添富逆向投资混合D",           "21.25", "92.44", "3.97", "0.8436", 10),
    @(2,  "470098", "汇添富逆向投资混合A",           "20.96", "92.44", "3.97", "0.8321", 10),
    @(3,  "001816", "汇添富新睿精选灵活配置混合A",    "1.64",  "94.71", "4.36", "0.0715", 10),
    @(4,  "001707", "诺安高端制造股票A",             "1.12",  "90.57", "6.18", "0.0692", 2),
    @(5,  "015784", "中信建投中证1000指数增强A",      "8.10",  "92.20", "0.68", "0.0551", 4),
    @(6,  "014509", "汇添富先进制造混合C",           "0.85",  "87.50", "5.35", "0.0455", 4),
    @(7,  "014508", "汇添富先进制造混合A",           "0.83",  "87.50", "5.35", "0.0444", 4),
    @(8,  "002164", "汇添富新睿精选灵活配置混合C",    "0.69",  "94.71", "4.36", "0.0301", 10),
    @(9,  "015785", "中信建投中证1000指数增强C",      "3.32",  "92.20", "0.68", "0.0226", 4),
    @(10, "015181", "汇添富逆向投资混合C",           "0.24",  "92.44", "3.97", "0.0095", 10),
    @(11, "014536", "诺安高端制造股票C",             "0.00",  "90.57", "6.18", "__NUM0__", 2)
)

$r = 2
foreach ($row in $rows) {
    $ws_new.Cells.Item($r, 1).Value = $row[0]
    $ws_new.Cells.Item($r, 2).Value = $row[1]
    $ws_new.Cells.Item($r, 3).Value = $row[2]
    $ws_new.Cells.Item($r, 4).Value = $row[3]
    $ws_new.Cells.Item($r, 5).Value = $row[4]
    $ws_new.Cells.Item($r, 6).Value = $row[5]
    if ($row[6] -eq "__NUM0__") {
        $ws_new.Cells.Item($r, 7).Value = 0
    } else {
        $ws_new.Cells.Item($r, 7).Value = $row[6]
    }
    $ws_new.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
